# Applies the "script re-run" update to the Ligue 1 2023-2024 odds sheet:
#  1) Several 2/3-row groups of match data (columns F:V) get re-ordered
#     in place (the underlying matches are the same set, only the row
#     order inside the group changed - columns A:E, which already agree
#     across every row of a group, are left untouched).
#  2) Eight brand-new match rows (107-114) are appended at the bottom.
#
# NOTE: custom functions with named parameters misbehave in this
# PowerShell host (parameters come through as $null), so everything
# below is written as flat, sequential statements using positional /
# direct-index access only.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Re-order the F:V block of a handful of rows. Read every row's
#    current F:V values first (so the reads never see a value already
#    overwritten by this same script), then write them back in the new
#    order.
# ---------------------------------------------------------------------

# --- group: rows 15,16,17 -> new15=old16, new16=old17, new17=old15 ---
$g1r15 = $ws.Range("F15:V15").Value()
$g1r16 = $ws.Range("F16:V16").Value()
$g1r17 = $ws.Range("F17:V17").Value()
$ws.Range("F15:V15").Value = $g1r16
$ws.Range("F16:V16").Value = $g1r17
$ws.Range("F17:V17").Value = $g1r15

# --- group: rows 42,43,44 -> new42=old43, new43=old44, new44=old42 ---
$g2r42 = $ws.Range("F42:V42").Value()
$g2r43 = $ws.Range("F43:V43").Value()
$g2r44 = $ws.Range("F44:V44").Value()
$ws.Range("F42:V42").Value = $g2r43
$ws.Range("F43:V43").Value = $g2r44
$ws.Range("F44:V44").Value = $g2r42

# --- group: rows 51,52 -> swap ---
$g3r51 = $ws.Range("F51:V51").Value()
$g3r52 = $ws.Range("F52:V52").Value()
$ws.Range("F51:V51").Value = $g3r52
$ws.Range("F52:V52").Value = $g3r51

# --- group: rows 60,61 -> swap ---
$g4r60 = $ws.Range("F60:V60").Value()
$g4r61 = $ws.Range("F61:V61").Value()
$ws.Range("F60:V60").Value = $g4r61
$ws.Range("F61:V61").Value = $g4r60

# --- group: rows 77,78,79 -> new77=old79, new78=old77, new79=old78 ---
$g5r77 = $ws.Range("F77:V77").Value()
$g5r78 = $ws.Range("F78:V78").Value()
$g5r79 = $ws.Range("F79:V79").Value()
$ws.Range("F77:V77").Value = $g5r79
$ws.Range("F78:V78").Value = $g5r77
$ws.Range("F79:V79").Value = $g5r78

# --- group: rows 86,87,88 -> new86=old87, new87=old88, new88=old86 ---
$g6r86 = $ws.Range("F86:V86").Value()
$g6r87 = $ws.Range("F87:V87").Value()
$g6r88 = $ws.Range("F88:V88").Value()
$ws.Range("F86:V86").Value = $g6r87
$ws.Range("F87:V87").Value = $g6r88
$ws.Range("F88:V88").Value = $g6r86

# --- group: rows 94,95,96 -> new94=old96, new95=old94, new96=old95 ---
$g7r94 = $ws.Range("F94:V94").Value()
$g7r95 = $ws.Range("F95:V95").Value()
$g7r96 = $ws.Range("F96:V96").Value()
$ws.Range("F94:V94").Value = $g7r96
$ws.Range("F95:V95").Value = $g7r94
$ws.Range("F96:V96").Value = $g7r95

# ---------------------------------------------------------------------
# 2) Append the new rows (107-114) with full A:V data.
# ---------------------------------------------------------------------

$newBlock = New-Object 'object[,]' 8,22

# Row 107 (Indice 106): PSG 5 - 2 Monaco
$newBlock[0,0]=106;  $newBlock[0,1]="france"; $newBlock[0,2]="ligue-1"; $newBlock[0,3]="2023-2024"; $newBlock[0,4]=45254.875
$newBlock[0,5]="PSG"; $newBlock[0,6]=5; $newBlock[0,7]="Monaco"; $newBlock[0,8]=2
$newBlock[0,9]=1.49;  $newBlock[0,10]="05/11/2023 11:03"; $newBlock[0,11]=1.48; $newBlock[0,12]="24/11/2023 20:58"
$newBlock[0,13]=4.92; $newBlock[0,14]="05/11/2023 11:03"; $newBlock[0,15]=5.17; $newBlock[0,16]="24/11/2023 20:59"
$newBlock[0,17]=5.91; $newBlock[0,18]="05/11/2023 11:03"; $newBlock[0,19]=6.11; $newBlock[0,20]="24/11/2023 20:59"
$newBlock[0,21]="https://www.betexplorer.com/football/france/ligue-1/psg-monaco/8AJ2LGCa/"

# Row 108 (Indice 107): Clermont 0 - 3 Lens
$newBlock[1,0]=107;  $newBlock[1,1]="france"; $newBlock[1,2]="ligue-1"; $newBlock[1,3]="2023-2024"; $newBlock[1,4]=45255.70833333334
$newBlock[1,5]="Clermont"; $newBlock[1,6]=0; $newBlock[1,7]="Lens"; $newBlock[1,8]=3
$newBlock[1,9]=3.44; $newBlock[1,10]="05/11/2023 11:03"; $newBlock[1,11]=4.03; $newBlock[1,12]="25/11/2023 16:55"
$newBlock[1,13]=3.4; $newBlock[1,14]="05/11/2023 11:03"; $newBlock[1,15]=3.51; $newBlock[1,16]="25/11/2023 16:55"
$newBlock[1,17]=2.2; $newBlock[1,18]="05/11/2023 11:03"; $newBlock[1,19]=2.02; $newBlock[1,20]="25/11/2023 16:55"
$newBlock[1,21]="https://www.betexplorer.com/football/france/ligue-1/clermont-lens/GGtUyxct/"

# Row 109 (Indice 108): Strasbourg 1 - 1 Marseille
$newBlock[2,0]=108; $newBlock[2,1]="france"; $newBlock[2,2]="ligue-1"; $newBlock[2,3]="2023-2024"; $newBlock[2,4]=45255.875
$newBlock[2,5]="Strasbourg"; $newBlock[2,6]=1; $newBlock[2,7]="Marseille"; $newBlock[2,8]=1
$newBlock[2,9]=4.03; $newBlock[2,10]="05/11/2023 11:03"; $newBlock[2,11]=4.28; $newBlock[2,12]="25/11/2023 20:38"
$newBlock[2,13]=3.66; $newBlock[2,14]="05/11/2023 11:03"; $newBlock[2,15]=3.67; $newBlock[2,16]="25/11/2023 20:38"
$newBlock[2,17]=1.93; $newBlock[2,18]="05/11/2023 11:03"; $newBlock[2,19]=1.91; $newBlock[2,20]="25/11/2023 20:35"
$newBlock[2,21]="https://www.betexplorer.com/football/france/ligue-1/strasbourg-marseille/jZ87KzS5/"

# Row 110 (Indice 109): Nice 1 - 0 Toulouse
$newBlock[3,0]=109; $newBlock[3,1]="france"; $newBlock[3,2]="ligue-1"; $newBlock[3,3]="2023-2024"; $newBlock[3,4]=45256.54166666666
$newBlock[3,5]="Nice"; $newBlock[3,6]=1; $newBlock[3,7]="Toulouse"; $newBlock[3,8]=0
$newBlock[3,9]=1.56; $newBlock[3,10]="05/11/2023 11:03"; $newBlock[3,11]=1.53; $newBlock[3,12]="26/11/2023 12:56"
$newBlock[3,13]=4.29; $newBlock[3,14]="05/11/2023 11:03"; $newBlock[3,15]=4.29; $newBlock[3,16]="26/11/2023 12:59"
$newBlock[3,17]=5.98; $newBlock[3,18]="05/11/2023 11:03"; $newBlock[3,19]=6.95; $newBlock[3,20]="26/11/2023 12:59"
$newBlock[3,21]="https://www.betexplorer.com/football/france/ligue-1/nice-toulouse/E7f2tIlP/"

# Row 111 (Indice 110): Montpellier 1 - 3 Brest
$newBlock[4,0]=110; $newBlock[4,1]="france"; $newBlock[4,2]="ligue-1"; $newBlock[4,3]="2023-2024"; $newBlock[4,4]=45256.625
$newBlock[4,5]="Montpellier"; $newBlock[4,6]=1; $newBlock[4,7]="Brest"; $newBlock[4,8]=3
$newBlock[4,9]=1.97; $newBlock[4,10]="05/11/2023 11:03"; $newBlock[4,11]=2.34; $newBlock[4,12]="26/11/2023 14:59"
$newBlock[4,13]=3.58; $newBlock[4,14]="05/11/2023 11:03"; $newBlock[4,15]=3.44; $newBlock[4,16]="26/11/2023 14:59"
$newBlock[4,17]=3.93; $newBlock[4,18]="05/11/2023 11:03"; $newBlock[4,19]=3.23; $newBlock[4,20]="26/11/2023 14:59"
$newBlock[4,21]="https://www.betexplorer.com/football/france/ligue-1/montpellier-brest/MDqxzGRh/"

# Row 112 (Indice 111): Nantes 0 - 0 Le Havre
$newBlock[5,0]=111; $newBlock[5,1]="france"; $newBlock[5,2]="ligue-1"; $newBlock[5,3]="2023-2024"; $newBlock[5,4]=45256.625
$newBlock[5,5]="Nantes"; $newBlock[5,6]=0; $newBlock[5,7]="Le Havre"; $newBlock[5,8]=0
$newBlock[5,9]=2.18; $newBlock[5,10]="05/11/2023 11:03"; $newBlock[5,11]=2; $newBlock[5,12]="26/11/2023 14:57"
$newBlock[5,13]=3.31; $newBlock[5,14]="05/11/2023 11:03"; $newBlock[5,15]=3.36; $newBlock[5,16]="26/11/2023 14:59"
$newBlock[5,17]=3.57; $newBlock[5,18]="05/11/2023 11:03"; $newBlock[5,19]=4.35; $newBlock[5,20]="26/11/2023 14:59"
$newBlock[5,21]="https://www.betexplorer.com/football/france/ligue-1/nantes-le-havre/AwtYzdCn/"

# Row 113 (Indice 112): Lorient 2 - 3 Metz
$newBlock[6,0]=112; $newBlock[6,1]="france"; $newBlock[6,2]="ligue-1"; $newBlock[6,3]="2023-2024"; $newBlock[6,4]=45256.625
$newBlock[6,5]="Lorient"; $newBlock[6,6]=2; $newBlock[6,7]="Metz"; $newBlock[6,8]=3
$newBlock[6,9]=1.88; $newBlock[6,10]="05/11/2023 11:03"; $newBlock[6,11]=2.1; $newBlock[6,12]="26/11/2023 14:59"
$newBlock[6,13]=3.63; $newBlock[6,14]="05/11/2023 11:03"; $newBlock[6,15]=3.35; $newBlock[6,16]="26/11/2023 14:57"
$newBlock[6,17]=4.25; $newBlock[6,18]="05/11/2023 11:03"; $newBlock[6,19]=3.94; $newBlock[6,20]="26/11/2023 14:59"
$newBlock[6,21]="https://www.betexplorer.com/football/france/ligue-1/lorient-metz/baCFIEdI/"

# Row 114 (Indice 113): Rennes 3 - 1 Reims
$newBlock[7,0]=113; $newBlock[7,1]="france"; $newBlock[7,2]="ligue-1"; $newBlock[7,3]="2023-2024"; $newBlock[7,4]=45256.71180555555
$newBlock[7,5]="Rennes"; $newBlock[7,6]=3; $newBlock[7,7]="Reims"; $newBlock[7,8]=1
$newBlock[7,9]=1.72; $newBlock[7,10]="05/11/2023 11:03"; $newBlock[7,11]=2.13; $newBlock[7,12]="26/11/2023 16:58"
$newBlock[7,13]=3.79; $newBlock[7,14]="05/11/2023 11:03"; $newBlock[7,15]=3.62; $newBlock[7,16]="26/11/2023 16:58"
$newBlock[7,17]=4.78; $newBlock[7,18]="05/11/2023 11:03"; $newBlock[7,19]=3.54; $newBlock[7,20]="26/11/2023 17:03"
$newBlock[7,21]="https://www.betexplorer.com/football/france/ligue-1/rennes-reims/4xDBJfsC/"

$ws.Range("A107:V114").Value = $newBlock

# Match the look & feel of the existing data: bold/centered/bordered
# index column (A) and a date-time formatted kickoff column (E).
$ws.Range("A106").Copy()
$ws.Range("A107:A114").PasteSpecial(-4122)

$ws.Range("E107:E114").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$excel.CutCopyMode = 0
